# Generate Report for Handoff
# Update Priority ("low" -> "ht") and Latest Handoff Datetime for the
# four handed-off files (rows 4-7) on both the zh-cn and de-de sheets.
# The Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
# de-de handoff datetime for these rows, so it is refreshed too.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

for ($r = 4; $r -le 7; $r++) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-09-01 16:38:16"

    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-09-01 16:38:22"

    $overview.Range("G$r").Value = "2016-09-01 16:38:22"
}
